$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price (column D) and Volume(1h) (column E) figures
# for the GitHub Actions symbol-list refresh.
$updates = @(
    @{ Cell = "D2"; Value = "308.06" },
    @{ Cell = "E2"; Value = "1.78%" },
    @{ Cell = "D3"; Value = "36.02" },
    @{ Cell = "E3"; Value = "3.05%" },
    @{ Cell = "D4"; Value = "5.120" },
    @{ Cell = "E4"; Value = "1.16%" },
    @{ Cell = "D5"; Value = "0.08101" },
    @{ Cell = "E5"; Value = "1.72%" },
    @{ Cell = "D6"; Value = "1.954" },
    @{ Cell = "E6"; Value = "1.01%" },
    @{ Cell = "D7"; Value = "7.758" },
    @{ Cell = "E7"; Value = "0.21%" },
    @{ Cell = "D8"; Value = "0.9303" },
    @{ Cell = "E8"; Value = "0.92%" },
    @{ Cell = "D9"; Value = "0.1386" },
    @{ Cell = "E9"; Value = "13.26%" },
    @{ Cell = "D10"; Value = "0.1915" },
    @{ Cell = "E10"; Value = "4.16%" },
    @{ Cell = "D11"; Value = "0.09219" },
    @{ Cell = "E11"; Value = "-1.34%" },
    @{ Cell = "D12"; Value = "0.03422" },
    @{ Cell = "E12"; Value = "-2.88%" },
    @{ Cell = "D13"; Value = "0.09838" },
    @{ Cell = "E13"; Value = "-0.07%" },
    @{ Cell = "D14"; Value = "0.001449" },
    @{ Cell = "E14"; Value = "4.17%" },
    @{ Cell = "D15"; Value = "0.005839" },
    @{ Cell = "E15"; Value = "0.16%" },
    @{ Cell = "D16"; Value = "3.620" },
    @{ Cell = "E16"; Value = "3.56%" },
    @{ Cell = "D17"; Value = "4.189" },
    @{ Cell = "E17"; Value = "3.51%" },
    @{ Cell = "E18"; Value = "1.47%" },
    @{ Cell = "D19"; Value = "0.3440" },
    @{ Cell = "E19"; Value = "-0.17%" },
    @{ Cell = "D20"; Value = "0.1342" },
    @{ Cell = "E20"; Value = "4.05%" },
    @{ Cell = "D21"; Value = "4.903" },
    @{ Cell = "E21"; Value = "-2.56%" },
    @{ Cell = "E22"; Value = "1.74%" },
    @{ Cell = "D23"; Value = "0.04448" },
    @{ Cell = "E23"; Value = "-1.06%" },
    @{ Cell = "D24"; Value = "0.001219" },
    @{ Cell = "E24"; Value = "0.24%" },
    @{ Cell = "D25"; Value = "0.004840" },
    @{ Cell = "D26"; Value = "0.0001243" },
    @{ Cell = "E26"; Value = "-0.52%" },
    @{ Cell = "D39"; Value = "0.02024" },
    @{ Cell = "E39"; Value = "5.70%" },
    @{ Cell = "D40"; Value = "0.04939" },
    @{ Cell = "E40"; Value = "4.07%" },
    @{ Cell = "D41"; Value = "0.007764" },
    @{ Cell = "E41"; Value = "2.76%" },
    @{ Cell = "D42"; Value = "0.01032" },
    @{ Cell = "E42"; Value = "8.07%" },
    @{ Cell = "E43"; Value = "3.88%" },
    @{ Cell = "D44"; Value = "0.002105" },
    @{ Cell = "E44"; Value = "-0.19%" },
    @{ Cell = "D45"; Value = "0.01128" },
    @{ Cell = "E45"; Value = "1.37%" },
    @{ Cell = "D46"; Value = "0.00006462" },
    @{ Cell = "E46"; Value = "3.08%" },
    @{ Cell = "E47"; Value = "0.27%" },
    @{ Cell = "D49"; Value = "0.001193" },
    @{ Cell = "E49"; Value = "-8.53%" },
    @{ Cell = "D50"; Value = "0.00002105" },
    @{ Cell = "E50"; Value = "0.27%" },
    @{ Cell = "D51"; Value = "0.0002005" },
    @{ Cell = "E51"; Value = "0.27%" }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    # Force text storage (matches source data which stores these as
    # plain text, e.g. "308.06" / "1.78%", not numeric/percentage values).
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
}
